$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (target stored OOXML width of 36; the runtime pads
# ColumnWidth by 5/6 when serializing, so back that out here)
$ws.Columns.Item(1).ColumnWidth = 35.166666666666664

# Replace the numeric header values with descriptive column names
$ws.Range("A1").Value = "aluno_id"
$ws.Range("B1").Value = "created_at"
$ws.Range("C1").Value = "nome"
$ws.Range("D1").Value = "updated_at"
$ws.Range("E1").Value = "curso_id"
$ws.Range("F1").Value = "matricula_id"
